$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the base task count (B3) per hito 3 progress update.
$ws.Range("B3").Value = 194

# Force a full recalculation so dependent formulas (B7, G7, I9, etc.)
# pick up the new value.
$excel.CalculateFullRebuild()
